# This script re-orders several data rows on the "Artfynd" sheet so that
# each row's full set of field values (Id, coordinates, times, comments,
# etc.) moves to a different row, matching the target revision.
#
# The sheet's used columns span A:AY (1..51).
#
# Columns Y (25, "Startdatum") and AA (27, "Slutdatum") hold plain text
# that looks like an ISO date ("2026-01-30"). Assigning such text through
# Range.Value2 makes Excel silently reinterpret it as a real date serial
# number (changing the cell's stored type/format). Since every row
# touched by this edit already shares the exact same Y/AA text, those two
# columns are simply left untouched during the row rewrites below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51                 # AY
$skipCols = @(25, 27)         # Y, AA

function Get-RowValues($rowNum) {
    return $ws.Range("A$($rowNum):AY$($rowNum)").Value2
}

function Set-RowValues($rowNum, $vals) {
    $startIdx = 1
    for ($c = 1; $c -le ($lastCol + 1); $c++) {
        $isSkip = ($c -le $lastCol) -and ($skipCols -contains $c)
        if ($isSkip -or $c -gt $lastCol) {
            if ($c -gt $startIdx) {
                $segLen = $c - $startIdx
                $seg = New-Object 'object[,]' 1, $segLen
                for ($i = 0; $i -lt $segLen; $i++) {
                    $seg[0, $i] = $vals[1, $startIdx + $i]
                }
                $fromAddr = $ws.Cells.Item($rowNum, $startIdx).Address($false, $false)
                $toAddr = $ws.Cells.Item($rowNum, $c - 1).Address($false, $false)
                $ws.Range("$fromAddr`:$toAddr").Value2 = $seg
            }
            $startIdx = $c + 1
        }
    }
}

# --- Rows 5 <-> 6 : simple swap ------------------------------------------
$r5 = Get-RowValues 5
$r6 = Get-RowValues 6
Set-RowValues 5 $r6
Set-RowValues 6 $r5

# --- Rows 13,14,15,16 : cyclic rotation (13->14->15->16->13) -------------
# i.e. new row14 = old row13, new row15 = old row14,
#      new row16 = old row15, new row13 = old row16
$r13 = Get-RowValues 13
$r14 = Get-RowValues 14
$r15 = Get-RowValues 15
$r16 = Get-RowValues 16
Set-RowValues 14 $r13
Set-RowValues 15 $r14
Set-RowValues 16 $r15
Set-RowValues 13 $r16

# --- Rows 17,18,19,20 : cyclic rotation (17->19->20->18->17) -------------
# i.e. new row19 = old row17, new row20 = old row19,
#      new row18 = old row20, new row17 = old row18
$r17 = Get-RowValues 17
$r18 = Get-RowValues 18
$r19 = Get-RowValues 19
$r20 = Get-RowValues 20
Set-RowValues 19 $r17
Set-RowValues 20 $r19
Set-RowValues 18 $r20
Set-RowValues 17 $r18

# --- Rows 25 <-> 26 : simple swap -----------------------------------------
$r25 = Get-RowValues 25
$r26 = Get-RowValues 26
Set-RowValues 25 $r26
Set-RowValues 26 $r25

# --- Rows 27 <-> 30 : simple swap -----------------------------------------
$r27 = Get-RowValues 27
$r30 = Get-RowValues 30
Set-RowValues 27 $r30
Set-RowValues 30 $r27

# --- Rows 32 <-> 34 : simple swap -----------------------------------------
$r32 = Get-RowValues 32
$r34 = Get-RowValues 34
Set-RowValues 32 $r34
Set-RowValues 34 $r32

Write-Host "Row re-ordering complete"
